$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$logMessages = @(
    "RandomPassword: Jz3lln4lRandomEmail: fPOojO@gmailAccountCreated",
    "RandomPassword: Jz3lln4lRandomEmail: fPOojO@gmailAccountDeleted",
    "RandomPassword: Ii5D9eczRandomEmail: SaBhRE@gmailAccountCreated",
    "RandomPassword: eeNScyHARandomEmail: stFpeZ@gmailAccountCreated",
    "RandomPassword: eeNScyHARandomEmail: stFpeZ@gmailAccountDeleted",
    "RandomPassword: sN2rpGmNRandomEmail: cKOBky@gmailAccountCreated",
    "RandomPassword: sN2rpGmNRandomEmail: cKOBky@gmailAccountDeleted",
    "RandomPassword: yDqo2EK5RandomEmail: PbQBjI@gmailAccountCreated",
    "RandomPassword: b8HuVog2RandomEmail: LzWZfY@gmailAccountCreated"
)

$startRow = 62
for ($i = 0; $i -lt $logMessages.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $logMessages[$i]
}
